$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0,1,1,0,1,0,1,0,1,0,0,1,1,0,0,0,1,0,0,0,0,1,1,1,0,0,1,0,0,0,0,1,0,0,0,1,0,1,1,1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
